$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Renewal transfer module update: the two "_Transfer" sheets drop their
#    free-text "Enter_Amount" column (the cash counterparts never had it).
# ---------------------------------------------------------------------------

# Jewel_GoldLoan_Renewal_Transfer: remove "Sanctioned_Amount" (col K) and
# "Enter_Amount" (col N). Delete the right-most column first so the other
# deletion still targets the intended column.
$wsJewelTransfer = $wb.Worksheets.Item("Jewel_GoldLoan_Renewal_Transfer")
$wsJewelTransfer.Range("N1").EntireColumn.Delete()
$wsJewelTransfer.Range("K1").EntireColumn.Delete()

# Jewel_Loan_Transcharge_Transfer: remove "Enter_Amount" (col I).
$wsTranscharge = $wb.Worksheets.Item("Jewel_Loan_Transcharge_Transfer")
$wsTranscharge.Range("I1").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2) Add the new "GL_Acc_Trans_Transfer" sheet (transfer counterpart of the
#    existing "GL_Acc_Trans_Cash" sheet), positioned right before it.
# ---------------------------------------------------------------------------
$wsGlCash = $wb.Worksheets.Item("GL_Acc_Trans_Cash")
$wsGlCash.Copy($wsGlCash)
$wsGlTransfer = $wb.Worksheets.Item("GL_Acc_Trans_Cash (2)")
$wsGlTransfer.Name = "GL_Acc_Trans_Transfer"

# Distinguish the new sheet's content from the Cash sheet it was copied from.
$wsGlTransfer.Range("A2").Value = "GL_Acc_Transaction_Transfer"

# Update the new sheet's selection.
$wsGlTransfer.Range("C6").Select()

# ---------------------------------------------------------------------------
# 3) Update the selections that moved on the edited "_Transfer" sheets.
# ---------------------------------------------------------------------------
$wsTranscharge.Range("H10").Select()

$wsJewelTransfer.Range("L11").Select()

# The edited Jewel_GoldLoan_Renewal_Transfer sheet becomes the active tab.
$wsJewelTransfer.Activate()
